$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.86"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.11"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.392"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05973"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.402"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.478"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8129"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9064"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1419"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07411"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03325"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03072"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.846"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001580"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04633"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006120"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005030"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009802"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00007795"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0002898"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.616"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03886"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006196"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002798"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007196"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005190"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002259"

$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
